$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 126, shifting existing rows 126-156 down to 127-157.
$ws.Rows.Item(126).Insert()

# Populate the newly inserted row 126 with the new data record.
$ws.Range("A126").Value = 4
$ws.Range("B126").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C126").Value = "Los Lagos"
$ws.Range("D126").Value = 44932
$ws.Range("E126").Value = 10
$ws.Range("F126").Value = 100112052
$ws.Range("G126").Value = "Albahaca"
$ws.Range("H126").Value = "Sin especificar"
$ws.Range("I126").Value = "Primera"
$ws.Range("J126").Value = 70
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = 6000
$ws.Range("N126").Value = "$/docena de matas"
$ws.Range("O126").Value = "Región Metropolitana"
$ws.Range("P126").Value = 1000
$ws.Range("Q126").Value = 6
$ws.Range("R126").Value = "Hortaliza"
